# Refresh market-price / profit columns (H:N) across all job sheets.
# Source data: scheduled market-board scrape; this mirrors a bulk value update
# with no formula or formatting changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 26.6
$ws.Range("I5").Value = 8.25
$ws.Range("K5").Value = 8.25
$ws.Range("M5").Value = 106.75
$ws.Range("H86").Value = 2636.2632
$ws.Range("I86").Value = 763.625
$ws.Range("J86").Value = 3998.182
$ws.Range("K86").Value = 763.625
$ws.Range("L86").Value = 3998.182
$ws.Range("M86").Value = 359.375
$ws.Range("N86").Value = -6244.182
$ws.Range("H89").Value = 2636.2632
$ws.Range("I89").Value = 763.625
$ws.Range("J89").Value = 3998.182
$ws.Range("K89").Value = 3818.125
$ws.Range("L89").Value = 19990.91
$ws.Range("M89").Value = 1797.875
$ws.Range("N89").Value = -31222.91
$ws.Range("H98").Value = 3256.2424
$ws.Range("I98").Value = 2889.25
$ws.Range("J98").Value = 15000
$ws.Range("K98").Value = 2889.25
$ws.Range("L98").Value = 15000
$ws.Range("M98").Value = -1391.25
$ws.Range("N98").Value = -17996
$ws.Range("H116").Value = 2342.5
$ws.Range("I116").Value = 2342.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2342.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1099.5
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 3256.2424
$ws.Range("I122").Value = 2889.25
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 8667.75
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -6217.75
$ws.Range("N122").Value = -49900
$ws.Range("H138").Value = 2381.2163
$ws.Range("I138").Value = 1973.5
$ws.Range("J138").Value = 2979.2
$ws.Range("K138").Value = 5920.5
$ws.Range("L138").Value = 8937.599999999999
$ws.Range("M138").Value = -780.5
$ws.Range("N138").Value = -19217.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2453553.8
$ws.Range("I2").Value = 2379.2856
$ws.Range("J2").Value = 5885198
$ws.Range("K2").Value = 2379.2856
$ws.Range("L2").Value = 5885198
$ws.Range("M2").Value = -2266.2856
$ws.Range("N2").Value = -5885424
$ws.Range("H5").Value = 119.5
$ws.Range("I5").Value = 122.4
$ws.Range("J5").Value = 117.42857
$ws.Range("K5").Value = 122.4
$ws.Range("L5").Value = 117.42857
$ws.Range("M5").Value = -10.40000000000001
$ws.Range("N5").Value = -341.42857
$ws.Range("H32").Value = 5010.67
$ws.Range("I32").Value = 4959.4536
$ws.Range("J32").Value = 6666.6665
$ws.Range("K32").Value = 4959.4536
$ws.Range("L32").Value = 6666.6665
$ws.Range("M32").Value = -4672.4536
$ws.Range("N32").Value = -7240.6665
$ws.Range("H61").Value = 5637.579
$ws.Range("I61").Value = 6960
$ws.Range("J61").Value = 4168.222
$ws.Range("K61").Value = 6960
$ws.Range("L61").Value = 4168.222
$ws.Range("M61").Value = -6748
$ws.Range("N61").Value = -4592.222
$ws.Range("H74").Value = 16388.8
$ws.Range("I74").Value = 1466.3334
$ws.Range("J74").Value = 22784.143
$ws.Range("K74").Value = 1466.3334
$ws.Range("L74").Value = 22784.143
$ws.Range("M74").Value = -592.3334
$ws.Range("N74").Value = -24532.143
$ws.Range("H77").Value = 16388.8
$ws.Range("I77").Value = 1466.3334
$ws.Range("J77").Value = 22784.143
$ws.Range("K77").Value = 7331.666999999999
$ws.Range("L77").Value = 113920.715
$ws.Range("M77").Value = -2963.666999999999
$ws.Range("N77").Value = -122656.715
$ws.Range("H116").Value = 2453553.8
$ws.Range("I116").Value = 2379.2856
$ws.Range("J116").Value = 5885198
$ws.Range("K116").Value = 2379.2856
$ws.Range("L116").Value = 5885198
$ws.Range("M116").Value = -85.28560000000016
$ws.Range("N116").Value = -5889786
$ws.Range("H136").Value = 5637.579
$ws.Range("I136").Value = 6960
$ws.Range("J136").Value = 4168.222
$ws.Range("K136").Value = 20880
$ws.Range("L136").Value = 12504.666
$ws.Range("M136").Value = -18330
$ws.Range("N136").Value = -17604.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2453553.8
$ws.Range("I3").Value = 2379.2856
$ws.Range("J3").Value = 5885198
$ws.Range("K3").Value = 2379.2856
$ws.Range("L3").Value = 5885198
$ws.Range("M3").Value = -2265.2856
$ws.Range("N3").Value = -5885426
$ws.Range("H4").Value = 119.5
$ws.Range("I4").Value = 122.4
$ws.Range("J4").Value = 117.42857
$ws.Range("K4").Value = 122.4
$ws.Range("L4").Value = 117.42857
$ws.Range("M4").Value = -7.400000000000006
$ws.Range("N4").Value = -347.42857
$ws.Range("H86").Value = 1369.2285
$ws.Range("I86").Value = 1323.5454
$ws.Range("J86").Value = 1446.5385
$ws.Range("K86").Value = 1323.5454
$ws.Range("L86").Value = 1446.5385
$ws.Range("M86").Value = -200.5454
$ws.Range("N86").Value = -3692.5385
$ws.Range("H89").Value = 1369.2285
$ws.Range("I89").Value = 1323.5454
$ws.Range("J89").Value = 1446.5385
$ws.Range("K89").Value = 6617.727
$ws.Range("L89").Value = 7232.692500000001
$ws.Range("M89").Value = -1001.727
$ws.Range("N89").Value = -18464.6925
$ws.Range("H134").Value = 159342.86
$ws.Range("I134").Value = 276750
$ws.Range("J134").Value = 2800
$ws.Range("K134").Value = 830250
$ws.Range("L134").Value = 8400
$ws.Range("M134").Value = -827715
$ws.Range("N134").Value = -13470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 997.5
$ws.Range("I2").Value = 997.5
$ws.Range("K2").Value = 997.5
$ws.Range("M2").Value = -884.5
$ws.Range("H7").Value = 21.8
$ws.Range("I7").Value = 15
$ws.Range("J7").Value = 49
$ws.Range("K7").Value = 15
$ws.Range("L7").Value = 49
$ws.Range("M7").Value = 98
$ws.Range("N7").Value = -275
$ws.Range("H31").Value = 2336.9756
$ws.Range("I31").Value = 1034.8077
$ws.Range("K31").Value = 1034.8077
$ws.Range("M31").Value = -739.8077000000001
$ws.Range("H34").Value = 2336.9756
$ws.Range("I34").Value = 1034.8077
$ws.Range("K34").Value = 1034.8077
$ws.Range("M34").Value = -832.8077000000001
$ws.Range("H58").Value = 3031.9692
$ws.Range("I58").Value = 1225.7333
$ws.Range("J58").Value = 4580.1714
$ws.Range("K58").Value = 1225.7333
$ws.Range("L58").Value = 4580.1714
$ws.Range("M58").Value = -1022.7333
$ws.Range("N58").Value = -4986.1714
$ws.Range("H99").Value = 117189.54
$ws.Range("I99").Value = 57405.555
$ws.Range("J99").Value = 251703.5
$ws.Range("K99").Value = 57405.555
$ws.Range("L99").Value = 251703.5
$ws.Range("M99").Value = -55907.555
$ws.Range("N99").Value = -254699.5
$ws.Range("H126").Value = 117189.54
$ws.Range("I126").Value = 57405.555
$ws.Range("J126").Value = 251703.5
$ws.Range("K126").Value = 172216.665
$ws.Range("L126").Value = 755110.5
$ws.Range("M126").Value = -169746.665
$ws.Range("N126").Value = -760050.5
$ws.Range("H132").Value = 1701.5349
$ws.Range("I132").Value = 1536.8334
$ws.Range("K132").Value = 4610.5002
$ws.Range("M132").Value = -2080.5002
$ws.Range("H136").Value = 3031.9692
$ws.Range("I136").Value = 1225.7333
$ws.Range("J136").Value = 4580.1714
$ws.Range("K136").Value = 3677.199900000001
$ws.Range("L136").Value = 13740.5142
$ws.Range("M136").Value = -1127.199900000001
$ws.Range("N136").Value = -18840.5142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 821.2
$ws.Range("I5").Value = 657.45
$ws.Range("J5").Value = 984.95
$ws.Range("K5").Value = 1972.35
$ws.Range("L5").Value = 2954.85
$ws.Range("M5").Value = -1860.35
$ws.Range("N5").Value = -3178.85
$ws.Range("H11").Value = 1567.7778
$ws.Range("I11").Value = 185
$ws.Range("J11").Value = 4333.3335
$ws.Range("K11").Value = 555
$ws.Range("L11").Value = 13000.0005
$ws.Range("M11").Value = -415
$ws.Range("N11").Value = -13280.0005
$ws.Range("H40").Value = 188.92308
$ws.Range("H113").Value = 3844.303
$ws.Range("I113").Value = 631.7059
$ws.Range("J113").Value = 7257.6875
$ws.Range("K113").Value = 1895.1177
$ws.Range("L113").Value = 21773.0625
$ws.Range("M113").Value = 274.8822999999998
$ws.Range("N113").Value = -26113.0625
$ws.Range("H121").Value = 20000562
$ws.Range("I121").Value = 564.3333
$ws.Range("J121").Value = 50000556
$ws.Range("K121").Value = 1692.9999
$ws.Range("L121").Value = 150001668
$ws.Range("M121").Value = -382.9999
$ws.Range("N121").Value = -150004288
$ws.Range("H131").Value = 2740.8928
$ws.Range("J131").Value = 1756.7925
$ws.Range("L131").Value = 5270.377500000001
$ws.Range("N131").Value = -15350.3775
$ws.Range("H135").Value = 821.2
$ws.Range("I135").Value = 657.45
$ws.Range("J135").Value = 984.95
$ws.Range("K135").Value = 5917.05
$ws.Range("L135").Value = 8864.550000000001
$ws.Range("M135").Value = -3382.05
$ws.Range("N135").Value = -13934.55

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1185.75
$ws.Range("I113").Value = 1156
$ws.Range("J113").Value = 1275
$ws.Range("K113").Value = 1156
$ws.Range("L113").Value = 1275
$ws.Range("M113").Value = 1014
$ws.Range("N113").Value = -5615
$ws.Range("H122").Value = 3499.2144
$ws.Range("I122").Value = 3998.4211
$ws.Range("J122").Value = 2445.3333
$ws.Range("K122").Value = 11995.2633
$ws.Range("L122").Value = 7335.999899999999
$ws.Range("M122").Value = -9545.263300000001
$ws.Range("N122").Value = -12235.9999
$ws.Range("H132").Value = 3362
$ws.Range("I132").Value = 2270.6667
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 6812.000100000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -4282.000100000001
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1889.4166
$ws.Range("I61").Value = 1776
$ws.Range("J61").Value = 2229.6667
$ws.Range("K61").Value = 1776
$ws.Range("L61").Value = 2229.6667
$ws.Range("M61").Value = -1574
$ws.Range("N61").Value = -2633.6667
$ws.Range("H113").Value = 1889.4166
$ws.Range("I113").Value = 1776
$ws.Range("J113").Value = 2229.6667
$ws.Range("K113").Value = 1776
$ws.Range("L113").Value = 2229.6667
$ws.Range("M113").Value = 394
$ws.Range("N113").Value = -6569.6667
$ws.Range("H122").Value = 1850.4783
$ws.Range("I122").Value = 1672.4117
$ws.Range("K122").Value = 5017.2351
$ws.Range("M122").Value = -2567.2351

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3931.2632
$ws.Range("I132").Value = 5487.6855
$ws.Range("J132").Value = 1455.1364
$ws.Range("K132").Value = 16463.0565
$ws.Range("L132").Value = 4365.4092
$ws.Range("M132").Value = -13933.0565
$ws.Range("N132").Value = -9425.4092
$ws.Range("H136").Value = 12093.708
$ws.Range("I136").Value = 13362.2
$ws.Range("J136").Value = 5751.25
$ws.Range("K136").Value = 40086.60000000001
$ws.Range("L136").Value = 17253.75
$ws.Range("M136").Value = -37536.60000000001
$ws.Range("N136").Value = -22353.75

